# edit.ps1 - Applies the "neutral voter language" + Mautinoa reorder edit
# described by the target diff, against $word.ActiveDocument (Word COM OM).

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Professional summary paragraph:
#    "...errors affecting all Black and Asian-American voters, developed..."
#    -> "...errors affecting 50M voters, developed..."
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters,",
    2) | Out-Null

# -----------------------------------------------------------------
# 2) Siege Analytics bullet:
#    "...errors affecting all Black and Asian-American voters, developed..."
#    -> "...errors affecting " + bold/colored "50M" + " voters, developed..."
#    Scope the search to just after "Discovered systematic race coding
#    errors affecting" so we only touch this bullet (not the summary or
#    the later Impact line, which share similar wording).
# -----------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute(
    "Discovered systematic race coding errors affecting",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAnchor = $anchor.End

# Replace "all Black and Asian-American voters" with plain "50M voters" first
$sub = $d.Range($afterAnchor, $afterAnchor + 120)
$sub.Find.Execute(
    "all Black and Asian-American voters",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "50M voters", 2) | Out-Null

# Now re-find just "50M" in that same neighborhood and make it bold + colored;
# Word will split the run automatically so only "50M" carries the new rPr.
$anchor2 = $d.Content
$anchor2.Find.Execute(
    "Discovered systematic race coding errors affecting",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAnchor2 = $anchor2.End
$sub2 = $d.Range($afterAnchor2, $afterAnchor2 + 60)
$sub2.Find.Execute(
    "50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sub2.Font.Bold = 1
$sub2.Font.Color = 5258796   # RGB(0x2C,0x3E,0x50) == w:color 2C3E50

# -----------------------------------------------------------------
# 3) Move the "Software Engineer - Mautinoa Technologies" job block
#    (Heading3 + 4 paragraphs) from its current spot (after Salsa Labs,
#    right before "KEY PROJECTS") to right after the Siege Analytics
#    bullets (right before "Senior Analyst - Myers Research").
# -----------------------------------------------------------------
$findBlock = $d.Content
$findBlock.Find.Execute(
    "Software Engineer - Mautinoa Technologies",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$jobHeading = $findBlock.Paragraphs(1)
$p2 = $jobHeading.Next()
$p3 = $p2.Next()
$p4 = $p3.Next()
$p5 = $p4.Next()
$blockEndMarker = $p5.Next()   # paragraph right after the block (Heading2 "KEY PROJECTS")

$fullBlock = $d.Range($jobHeading.Range.Start, $blockEndMarker.Range.Start)
$fullBlock.Cut() | Out-Null

# Re-find the insertion anchor after the cut (positions shifted)
$destFind = $d.Content
$destFind.Find.Execute(
    "Senior Analyst - Myers Research",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$destPara = $destFind.Paragraphs(1)
$insertPoint = $d.Range($destPara.Range.Start, $destPara.Range.Start)
$insertPoint.Paste() | Out-Null

# Pasting at a collapsed insertion point drops the source paragraph's
# Heading3 style from the first pasted paragraph - restore it explicitly.
$fixFind = $d.Content
$fixFind.Find.Execute(
    "Software Engineer - Mautinoa Technologies",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fixFind.Paragraphs(1).set_Style("Heading 3") | Out-Null

# -----------------------------------------------------------------
# 4) Project impact line:
#    "Impact: Corrected demographic data affecting all Black and
#    Asian-American voters, improved electoral prediction accuracy by 22%"
#    -> "Impact: Corrected demographic data affecting 50M voters
#    nationwide, improved electoral prediction accuracy by 22%"
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2) | Out-Null

Write-Output "edit complete"
